$p = $ppt.ActivePresentation

# --- 1) Fix typo on the "12-Factor" matching slide (slide 13): ---
#     "does" -> "do<cross mark emoji>es", splitting the run in three.
$slide13 = $p.Slides.Item(13)
$contentShape = $slide13.Shapes.Item(2)
$bodyRange = $contentShape.TextFrame.TextRange
$crossMark = [string][char]0x274C
$wordRange = $bodyRange.Find("does")
$wordRange.Text = "do" + $crossMark + "es"

# --- 2) Append a new "Questions?" slide at the end of the deck ---
$newSlide = $p.Slides.Add($p.Slides.Count + 1, 2)

$titleShape = $newSlide.Shapes.Item(1)
$titleShape.TextFrame.TextRange.Text = "Questions?"

$bodyShape = $newSlide.Shapes.Item(2)
$bodyTextRange = $bodyShape.TextFrame.TextRange
$bodyTextRange.Text = "???"
$bodyTextRange.Font.Size = 140
$bodyTextRange.Font.Bold = $true
$bodyTextRange.ParagraphFormat.Alignment = 2
